$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Update price (D) and volume-1h (E) columns for each coin row,
# and swap the THORChain/LidoDAOToken rows (36/37) contents.
Set-TextCell "D2" "38.441.20"
Set-TextCell "E2" "  +1.95%  "
Set-TextCell "D3" "2.099.79"
Set-TextCell "E3" "  +3.69%  "
Set-TextCell "E4" "  -0.09%  "
Set-TextCell "D5" "228.93"
Set-TextCell "E5" "  +0.70%  "
Set-TextCell "D6" "0.615"
Set-TextCell "E6" "  +1.27%  "
Set-TextCell "D7" "61.44"
Set-TextCell "E7" "  +3.03%  "
Set-TextCell "D8" "0.999"
Set-TextCell "E8" "  -0.14%  "
Set-TextCell "D9" "0.381"
Set-TextCell "E9" "  +1.77%  "
Set-TextCell "D10" "0.0845"
Set-TextCell "E10" "  +2.98%  "
Set-TextCell "E11" "  +0.37%  "
Set-TextCell "D12" "2.410.65"
Set-TextCell "E12" "  +3.63%  "
Set-TextCell "D13" "14.81"
Set-TextCell "E13" "  +3.11%  "
Set-TextCell "D14" "22.33"
Set-TextCell "E14" "  +6.37%  "
Set-TextCell "D15" "0.781"
Set-TextCell "E15" "  +2.30%  "
Set-TextCell "D16" "5.48"
Set-TextCell "E16" "  +5.93%  "
Set-TextCell "D17" "2.100.76"
Set-TextCell "E17" "  +3.77%  "
Set-TextCell "D18" "38.353.47"
Set-TextCell "E18" "  +1.90%  "
Set-TextCell "D19" "6.02"
Set-TextCell "E19" "  +2.55%  "
Set-TextCell "D20" "70.41"
Set-TextCell "E20" "  +1.54%  "
Set-TextCell "D21" "0.0₃0835"
Set-TextCell "E21" "  +1.76%  "
Set-TextCell "D22" "225.70"
Set-TextCell "E22" "  +0.97%  "
Set-TextCell "E23" "  -0.02%  "
Set-TextCell "E24" "  +0.18%  "
Set-TextCell "D25" "2.31"
Set-TextCell "E25" "  +3.26%  "
Set-TextCell "D26" "169.83"
Set-TextCell "E26" "  +1.49%  "
Set-TextCell "D27" "9.42"
Set-TextCell "E27" "  +1.41%  "
Set-TextCell "E28" "  +0.86%  "
Set-TextCell "D29" "19.05"
Set-TextCell "E29" "  +1.59%  "
Set-TextCell "E30" "  +9.03%  "
Set-TextCell "E31" "  -0.05%  "
Set-TextCell "D32" "2.35"
Set-TextCell "E32" "  +6.56%  "
Set-TextCell "D33" "4.75"
Set-TextCell "E33" "  +6.39%  "
Set-TextCell "D34" "4.46"
Set-TextCell "E34" "  +2.21%  "
Set-TextCell "D35" "0.0605"
Set-TextCell "E35" "  +0.53%  "
Set-TextCell "B36" "LidoDAOToken"
Set-TextCell "C36" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D36" "2.39"
Set-TextCell "E36" "  +4.33%  "
Set-TextCell "B37" "THORChain"
Set-TextCell "C37" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell "D37" "6.41"
Set-TextCell "E37" "  +0.31%  "
Set-TextCell "E38" "  +3.10%  "
Set-TextCell "E39" "  -0.06%  "
Set-TextCell "D40" "18.21"
Set-TextCell "E40" "  +2.27%  "
Set-TextCell "D41" "1.536.26"
Set-TextCell "E41" "  +0.42%  "
Set-TextCell "D42" "99.93"
Set-TextCell "E42" "  +4.87%  "
Set-TextCell "E43" "  +2.17%  "
Set-TextCell "E44" "  +0.91%  "
Set-TextCell "D45" "0.0910"
Set-TextCell "E45" "  +0.44%  "
Set-TextCell "D46" "4.16"
Set-TextCell "E46" "  +2.71%  "
Set-TextCell "E47" "  +0.99%  "
Set-TextCell "D48" "7.48"
Set-TextCell "E48" "  +5.58%  "
Set-TextCell "D49" "1.04"
Set-TextCell "E49" "  +4.17%  "
Set-TextCell "E50" "  +0.80%  "
Set-TextCell "D51" "2.296.15"
Set-TextCell "E51" "  +3.61%  "

Write-Output "Updated cryptos list."
